$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "SAN DIEGO AREA TOTALS" label from B2 into A2, and set B2 to "Totals"
$ws.Range("A2").Value = $ws.Range("B2").Value2
$ws.Range("B2").Value = "Totals"

# Column A should now match column B's (auto-fit) width
$ws.Columns("A:B").EntireColumn.AutoFit()

# Update the selection to match the diff (column A selected)
$ws.Range("A1:A1048576").Select()
